$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52:131 down to 53:132.
# Excel automatically carries the formatting (e.g. the date style in column D)
# from the row above into the newly inserted row.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly price record.
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 44579
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 100112003
$ws.Range("G52").Value = "Ajo"
$ws.Range("H52").Value = "Chino"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 400
$ws.Range("K52").Value = 17000
$ws.Range("L52").Value = 18000
$ws.Range("M52").Value = 17500
$ws.Range("N52").Value = "$/caja 10 kilos"
$ws.Range("O52").Value = "China"
$ws.Range("P52").Value = 1750
$ws.Range("Q52").Value = 10
$ws.Range("R52").Value = "Hortaliza"
